$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autoavaliação")

# Fill the "Completo" data-validation selection into the evaluation columns
# B37:B53 and F37:F56 on the "Autoavaliação" sheet.
$ws.Range("B37:B53").Value = "Completo"
$ws.Range("F37:F56").Value = "Completo"

# Update view state to match target workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("J52").Select()
